$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "دی 98" section: remove the "* Meetings @Parsiss" activity entry (row 77) ---
# Keep the cell formatting/styles intact, just clear the contents so the row reads blank.
$ws.Range("B77").ClearContents()
$ws.Range("C77").ClearContents()

# --- "بهمن 98" section: add a new "* Registration" activity line ---
# Insert a new row before the existing "Total Hours" row (row 87), shifting the
# Total/Paid/Not Paid rows down by one.
$ws.Rows("87:87").Insert()

# Copy the formatting from the row above (Tracker DataCapture, row 86) into the
# newly inserted row so the new cells keep the same look (borders/alignment).
$ws.Range("B86:C86").Copy()
$ws.Range("B87:C87").PasteSpecial(-4122)

# Fill in the new activity row values.
$ws.Range("B87").Value2 = "* Registration"
$ws.Range("C87").Value2 = 1

# The "Total Hours" row (now shifted to row 88) must sum the extended range.
$ws.Range("C88").Formula = "=SUM(C84:C87)"

# The "Not Paid" hours total (now shifted to row 90) grows by the new entry's hour.
$ws.Range("D90").Value2 = 7

# Update the visible selection to mirror the authored view state.
$ws.Range("D91").Select()
